$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Tipo" column with values for each certificate row
$ws.Range("G1").Value = "Tipo"
$ws.Range("G2").Value = "Primario"
$ws.Range("G3").Value = "Secundario"
$ws.Range("G4").Value = "Vencido"

# Update the active selection to match the saved view state
$ws.Range("H7").Select()
